$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values that look numeric (e.g. "1.003")
# as well as multi-dot strings (e.g. "30.320.67") that must stay as literal
# text, matching the original workbook formatting. Force a text number
# format on each Price cell before writing so Excel does not reinterpret
# the string as a number.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.320.67'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.929.50'
$ws.Range("E3").Value = '  -0.04%  '

# Row 4
$ws.Range("E4").Value = '  +0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7431'
$ws.Range("E5").Value = '  +3.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.72'
$ws.Range("E6").Value = '  -2.29%  '

# Row 8
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.45'
$ws.Range("E8").Value = '  -2.06%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3152'
$ws.Range("E9").Value = '  -1.81%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07024'
$ws.Range("E10").Value = '  -1.26%  '

# Row 11
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08033'
$ws.Range("E11").Value = '  +0.32%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7761'
$ws.Range("E12").Value = '  -1.70%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.896.44'
$ws.Range("E13").Value = '  -1.70%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.379'
$ws.Range("E14").Value = '  -0.09%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.09'
$ws.Range("E15").Value = '  -1.87%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.50'
$ws.Range("E16").Value = '  -1.28%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.330.14'
$ws.Range("E17").Value = '  +0.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.992'
$ws.Range("E18").Value = '  +4.01%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '250.65'
$ws.Range("E19").Value = '  -2.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007920'
$ws.Range("E20").Value = '  -2.20%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.24%  '

# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.151.69'
$ws.Range("E22").Value = '  -1.37%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.645'
$ws.Range("E24").Value = '  -2.72%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.556'
$ws.Range("E25").Value = '  -0.13%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.47'
$ws.Range("E26").Value = '  +0.45%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.00'
$ws.Range("E27").Value = '  -0.68%  '

# Row 28
$ws.Range("E28").Value = '  +0.56%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.164'
$ws.Range("E29").Value = '  -5.46%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.565'
$ws.Range("E30").Value = '  +2.11%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.357'
$ws.Range("E31").Value = '  +0.18%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.406'
$ws.Range("E32").Value = '  -0.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.103'
$ws.Range("E33").Value = '  -1.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05214'
$ws.Range("E34").Value = '  +1.67%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.307'
$ws.Range("E35").Value = '  +1.91%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7536'
$ws.Range("E36").Value = '  +0.56%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.768'
$ws.Range("E37").Value = '  -0.19%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01949'
$ws.Range("E38").Value = '  -1.99%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.792'
$ws.Range("E39").Value = '  -0.28%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.504'
$ws.Range("E40").Value = '  +1.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.51'
$ws.Range("E41").Value = '  -2.31%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4496'
$ws.Range("E42").Value = '  -0.54%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.955'
$ws.Range("E43").Value = '  -2.18%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8424'
$ws.Range("E44").Value = '  -0.38%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.003'
$ws.Range("E45").Value = '  +0.28%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.673'
$ws.Range("E46").Value = '  +2.46%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.941'
$ws.Range("E47").Value = '  +1.51%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.22'
$ws.Range("E48").Value = '  -0.14%  '

# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.088.95'
$ws.Range("E49").Value = '  -0.22%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.53'
$ws.Range("E50").Value = '  +1.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1224'
$ws.Range("E51").Value = '  +6.78%  '
